$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2 value: 102 -> 104
$ws.Range("B2").Value = 104

# Update B11 value: 1.8 -> 5
$ws.Range("B11").Value = 5

# Update the active selection from B2 to B11
$ws.Range("B11").Select()
